$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.957.84'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = "'311.67"
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("D7").Value = "'0.5089"
$ws.Range("D8").Value = "'0.3804"
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("D9").Value = "'0.08247"
$ws.Range("E9").Value = '  -7.23%  '
$ws.Range("D10").Value = "'1.107"
$ws.Range("E10").Value = '  -0.97%  '
$ws.Range("B11").Value = 'Polkadot'
$ws.Range("C11").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D11").Value = "'6.196"
$ws.Range("E11").Value = '  -2.69%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = "'20.49"
$ws.Range("E12").Value = '  -0.91%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.854.24'
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = "'7.175"
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").Value = "'1.005"
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = "'0.00001096"
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = "'90.56"
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = "'0.06606"
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = "'17.62"
$ws.Range("E19").Value = '  -2.08%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = "'6.019"
$ws.Range("E21").Value = '  -1.45%  '
$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").Value = '27.977.44'
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = "'11.02"
$ws.Range("E23").Value = '  -3.95%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = "'2.242"
$ws.Range("E24").Value = '  -1.57%  '
$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D25").Value = "'2.535"
$ws.Range("E25").Value = '  +1.29%  '
$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.069.30'
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = "'157.86"
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'20.38"
$ws.Range("E28").Value = '  -1.26%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = "'124.38"
$ws.Range("E29").Value = '  -1.33%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = "'0.1052"
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = "'1.035"
$ws.Range("E31").Value = '  -1.63%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'5.613"
$ws.Range("E32").Value = '  +0.49%  '
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = "'3.595"
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("B34").Value = 'FraxShare'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D34").Value = "'9.388"
$ws.Range("E34").Value = '  +0.34%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = "'0.06514"
$ws.Range("E35").Value = '  -0.36%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = "'0.02403"
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = "'0.2161"
$ws.Range("E37").Value = '  -1.02%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = "'1.205"
$ws.Range("E38").Value = '  +0.61%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = "'0.6450"
$ws.Range("E39").Value = '  +1.27%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = "'1.224"
$ws.Range("E40").Value = '  -4.26%  '
$ws.Range("B41").Value = 'InternetComputer(DFINITY)'
$ws.Range("C41").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D41").Value = "'4.867"
$ws.Range("E41").Value = '  -0.86%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = "'11.13"
$ws.Range("E42").Value = '  -3.53%  '
$ws.Range("B43").Value = 'Decentraland'
$ws.Range("C43").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D43").Value = "'0.6075"
$ws.Range("E43").Value = '  +1.33%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'13.08"
$ws.Range("E44").Value = '  -0.61%  '
$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = "'1.278"
$ws.Range("E45").Value = '  -0.37%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = "'3.651"
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = "'1.996"
$ws.Range("E47").Value = '  +0.31%  '
$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D48").Value = "'1.206"
$ws.Range("E48").Value = '  -1.48%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = "'119.88"
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = "'78.56"
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.06828"
$ws.Range("E51").Value = '  -0.83%  '
